$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column M (inherits formatting from the adjacent 2020 column L)
# so the new cells automatically pick up the same styles used by column L.
$ws.Range("M1").EntireColumn.Insert()

# Populate the new "2021" data column
$ws.Range("M4").Value = 2021
$ws.Range("M5").Value = 98
$ws.Range("M6").Value = 97
$ws.Range("M7").Value = 96

# Reset the view back to the sheet's default top-left cell / selection
$ws.Range("A1").Select()
